# Calendario-Responsabilidades.xlsx — add "VerReportesController" progress row,
# mark ReportesController's "25/03/2024 al 31/03/2024" week complete, and
# refresh the sheet view (matches "Interfaz de reportes terminada (100%)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) ReportesController ("25/03/2024" al "31/03/2024") column now at 100%.
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = 100

# ---------------------------------------------------------------------------
# 2) New row for the finished reports screen controller.
#    Clone row 14's look (fonts/fills/borders/alignment) down into row 15,
#    then overwrite the values for the new controller.
# ---------------------------------------------------------------------------
$ws.Range("A14:E14").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A15").Value = "VerReportesController"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 85
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = 0

$ws.Rows.Item(15).RowHeight = 24

# ---------------------------------------------------------------------------
# 3) Row heights on the existing rows nudged slightly (14.4->15 default,
#    25.05->25.15 per data row) and column widths re-measured by Excel.
# ---------------------------------------------------------------------------
$ws.Range("A1:E1").RowHeight = 25.15
$ws.Range("A2:E2").RowHeight = 25.15
$ws.Range("A3:E3").RowHeight = 25.15
$ws.Range("A4:E4").RowHeight = 25.15
$ws.Range("A5:E5").RowHeight = 25.15
$ws.Range("A6:E6").RowHeight = 25.15
$ws.Range("A7:E7").RowHeight = 25.15
$ws.Range("A8:E8").RowHeight = 25.15
$ws.Range("A9:E9").RowHeight = 25.15
$ws.Range("A10:E10").RowHeight = 25.15
$ws.Range("A11:E11").RowHeight = 25.15
$ws.Range("A12:E12").RowHeight = 25.15
$ws.Range("A13:E13").RowHeight = 25.15
$ws.Range("A14:E14").RowHeight = 25.15

$ws.Columns.Item(1).ColumnWidth = 33.85546875
$ws.Columns.Item(3).ColumnWidth = 29.7109375
$ws.Columns.Item(4).ColumnWidth = 29.85546875
$ws.Columns.Item(5).ColumnWidth = 30.28515625

# ---------------------------------------------------------------------------
# 4) Conditional formatting: extend the colour scales/banding that used to
#    stop at row 14 so the new row 15 is covered too, mirroring what Excel
#    does when the last row of a CF range is duplicated downward.
# ---------------------------------------------------------------------------

# Column B colour scale B2:B14 -> B2:B15
$ws.Range("B2:B14").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B2:B15"))

# Red/Yellow/Green banding B2:E13 -> B2:E14 (absorbs the old standalone B14 rule)
$ws.Range("B2:E13").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B2:E14"))

# E column standalone colour scale E14 -> E14:E15
$ws.Range("E14").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E14:E15"))

# Drop the old standalone B14/C14(cellIs)/D14(cellIs)/E14(cellIs) rules; they
# get rebuilt fresh on row 15 below (C14/D14 colour-scale rules stay put).
$ws.Range("B14").FormatConditions.Delete()

$ws.Range("C14").FormatConditions.Item(2).Delete()

$ws.Range("D14").FormatConditions.Item(2).Delete()

$ws.Range("E14").FormatConditions.Item(2).Delete()

# Rebuild the per-cell red/yellow/green banding on row 15.
$b15 = $ws.Range("B15").FormatConditions
$r = $b15.Add(1, 1, "80", "100")
$r.Interior.Color = 6736998
$r = $b15.Add(1, 1, "0", "40")
$r.Interior.Color = 255
$r = $b15.Add(1, 1, "50", "70")
$r.Interior.Color = 5296274

$c15 = $ws.Range("C15").FormatConditions
$r = $c15.Add(1, 1, "80", "100")
$r.Interior.Color = 6736998
$r = $c15.Add(1, 1, "0", "40")
$r.Interior.Color = 255
$r = $c15.Add(1, 1, "50", "70")
$r.Interior.Color = 5296274

$d15 = $ws.Range("D15").FormatConditions
$r = $d15.Add(1, 1, "80", "100")
$r.Interior.Color = 6736998
$r = $d15.Add(1, 1, "0", "40")
$r.Interior.Color = 255
$r = $d15.Add(1, 1, "50", "70")
$r.Interior.Color = 5296274

$e15 = $ws.Range("E15").FormatConditions
$r = $e15.Add(1, 1, "80", "100")
$r.Interior.Color = 6736998
$r = $e15.Add(1, 1, "0", "40")
$r.Interior.Color = 255
$r = $e15.Add(1, 1, "50", "70")
$r.Interior.Color = 5296274

# C15/D15 also get their own colour-scale rule (mirrors C14/D14).
$c15cs = $ws.Range("C15").FormatConditions.AddColorScale(3)
$d15cs = $ws.Range("D15").FormatConditions.AddColorScale(3)

# ---------------------------------------------------------------------------
# 5) Sheet dimension / selection follow the new last row.
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
$excel.ActiveWindow.ScrollRow = 8
